$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 4 (rank 3): Nickname "Charo" -> "Esther SV", Points 2000 -> 4000
$ws.Range("B4").Value = "Esther SV"
$ws.Range("D4").Value = 4000

# Update row 6 (rank 5): Nickname "Guzm4n" -> "Joh14"
$ws.Range("B6").Value = "Joh14"

# Update selection to match the saved view state
$ws.Range("B7").Select()
